$wb = $excel.ActiveWorkbook

$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and worksheet name
# lookup by Item(name) is case-insensitive, so both names would resolve to
# the same (first) sheet. Use the 1-based tab index instead to disambiguate.
$wsVecbf = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF = $wb.Worksheets.Item(6)   # Vector_BF
$wsVecAlpha = $wb.Worksheets.Item("Vector_Alpha")

if ($wsVecbf.Name -ne "Vector_bf") { throw "expected Vector_bf at index 5, got $($wsVecbf.Name)" }
if ($wsVecBF.Name -ne "Vector_BF") { throw "expected Vector_BF at index 6, got $($wsVecBF.Name)" }

# All the numeric-looking values in this workbook are stored as plain TEXT
# (shared-string) cells, not numbers -- e.g. "4.5 - x" sits next to "-5.0",
# which is text too. Writing a numeric-looking string via .Value normally
# gets auto-converted to a real number by Excel, so we briefly force Text
# number-format, assign the value, then restore the original (default/
# General, unstyled) cell formatting by pasting formats from an untouched
# neighboring cell -- this keeps the cell's `s` (style) attribute identical
# to the rest of the sheet instead of leaving a "@"-formatted style behind.
function Set-TextValue($range, $text, $fmtSrc) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $fmtSrc.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# Restricciones_del_lider (fmt source: C1, untouched default-style text cell)
$liderFmt = $wsLider.Range("C1")
Set-TextValue $wsLider.Range("A2") "2.8499999999999996 - x" $liderFmt
Set-TextValue $wsLider.Range("B2") "-3.3499999999999996" $liderFmt
Set-TextValue $wsLider.Range("D2") "0.3" $liderFmt
Set-TextValue $wsLider.Range("A3") "-2.8499999999999996 + x" $liderFmt
Set-TextValue $wsLider.Range("B3") "2.3499999999999996" $liderFmt
Set-TextValue $wsLider.Range("D3") "0.09" $liderFmt

# Restricciones_del_follower (fmt source: C1, untouched default-style text cell)
$followerFmt = $wsFollower.Range("C1")
Set-TextValue $wsFollower.Range("A2") "4.440892098500626e-16 - 1.1102230246251565e-16y" $followerFmt
Set-TextValue $wsFollower.Range("B2") "-1.0000000000000004" $followerFmt
Set-TextValue $wsFollower.Range("D2") "0.19" $followerFmt
Set-TextValue $wsFollower.Range("E2") "0" $followerFmt
Set-TextValue $wsFollower.Range("F2") "0.0" $followerFmt
Set-TextValue $wsFollower.Range("A3") "-4.440892098500626e-16 + 1.1102230246251565e-16y" $followerFmt
Set-TextValue $wsFollower.Range("B3") "-0.9999999999999996" $followerFmt
Set-TextValue $wsFollower.Range("D3") "0.79" $followerFmt
Set-TextValue $wsFollower.Range("E3") "0" $followerFmt
Set-TextValue $wsFollower.Range("F3") "1.7999999999999998" $followerFmt

# Punto_modificado (fmt source: A1, untouched default-style text cell)
$puntoFmt = $wsPunto.Range("A1")
Set-TextValue $wsPunto.Range("A2") "2.8499999999999996" $puntoFmt
Set-TextValue $wsPunto.Range("B2") "4.449999999999999" $puntoFmt

# Vector_bf (fmt source: A1)
$vecbfFmt = $wsVecbf.Range("A1")
Set-TextValue $wsVecbf.Range("A2") "-2.8499999999999996" $vecbfFmt

# Vector_BF (fmt source: A1)
$vecBFFmt = $wsVecBF.Range("A1")
Set-TextValue $wsVecBF.Range("A2") "0.49" $vecBFFmt
Set-TextValue $wsVecBF.Range("A3") "1.0" $vecBFFmt

# Vector_Alpha: A2 is a genuine number both before and after, just a new value.
$wsVecAlpha.Range("A2").Value = 0.72
